# "se agrego la edad" - add the age under the (now capitalized + underlined)
# name, and leave a trailing blank underlined line.
#
# Before: a single paragraph containing the plain text "derekcyborg".
# After:  three paragraphs:
#   1) "Derekcyborg " - the name, capitalized, underlined, with the
#      spell-checker's wavy-underline markers around the word itself
#      (the trailing space is its own run, still underlined)
#   2) "Veinte años"  - the age, underlined
#   3) an empty underlined paragraph
#
# Word's object model doesn't give us a way to stamp <w:proofErr/> markers
# or to leave a completely run-less empty paragraph, so the new content is
# written with Range.InsertXML, which replaces the contents of the target
# Range with the supplied WordprocessingML - exactly what a "paste from
# Word" or "insert XML" COM call would do.

$d = $word.ActiveDocument

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:pPr><w:rPr><w:u w:val="single"/><w:lang w:val="es-MX"/></w:rPr></w:pPr>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:rPr><w:u w:val="single"/><w:lang w:val="es-MX"/></w:rPr><w:t>Derekcyborg</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:rPr><w:u w:val="single"/><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
          '</w:p>' +
          '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:pPr><w:rPr><w:u w:val="single"/><w:lang w:val="es-MX"/></w:rPr></w:pPr>' +
            '<w:r><w:rPr><w:u w:val="single"/><w:lang w:val="es-MX"/></w:rPr><w:t>Veinte años</w:t></w:r>' +
          '</w:p>' +
          '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:pPr><w:rPr><w:u w:val="single"/><w:lang w:val="es-MX"/></w:rPr></w:pPr>' +
          '</w:p>'

# Replace the whole (single) paragraph that holds "derekcyborg" with the
# three paragraphs above.
$d.Paragraphs(1).Range.InsertXML($newXml)
